$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header swap: average_doctor <-> average_doctor_old
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.201
$ws.Range("AJ4").Value = 0.065
$ws.Range("AK4").Value = 0.255
$ws.Range("AU4").Value = 0.147
$ws.Range("AV4").Value = 0.026
$ws.Range("AW4").Value = 0.163
$ws.Range("BA4").Value = 1.989
$ws.Range("BB4").Value = 0.166
$ws.Range("BC4").Value = 0.408
$ws.Range("BG4").Value = 0.745
$ws.Range("BH4").Value = 0.134
$ws.Range("BI4").Value = 0.366
$ws.Range("BM4").Value = 0.696
$ws.Range("BN4").Value = 0.083
$ws.Range("BO4").Value = 0.288
$ws.Range("BP4").Value = 0.663
$ws.Range("BQ4").Value = 0.664
$ws.Range("E4").Value = 0.382
$ws.Range("F4").Value = 0.08400000000000001
$ws.Range("G4").Value = 0.289
$ws.Range("N4").Value = 0.405
$ws.Range("O4").Value = 0.063
$ws.Range("P4").Value = 0.252
$ws.Range("W4").Value = 0.221
$ws.Range("X4").Value = 0.102
$ws.Range("Y4").Value = 0.32
$ws.Range("AI5").Value = 0.24
$ws.Range("AJ5").Value = 0.094
$ws.Range("AK5").Value = 0.307
$ws.Range("AU5").Value = 0.303
$ws.Range("AV5").Value = 0.1
$ws.Range("AW5").Value = 0.315
$ws.Range("BA5").Value = 1.394
$ws.Range("BC5").Value = 0.29
$ws.Range("BG5").Value = 0.418
$ws.Range("BH5").Value = 0.05
$ws.Range("BI5").Value = 0.223
$ws.Range("BM5").Value = 0.585
$ws.Range("BN5").Value = 0.073
$ws.Range("BO5").Value = 0.269
$ws.Range("BP5").Value = 0.465
$ws.Range("BQ5").Value = 0.461
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.101
$ws.Range("G5").Value = 0.318
$ws.Range("N5").Value = 0.764
$ws.Range("O5").Value = 0.077
$ws.Range("P5").Value = 0.278
$ws.Range("W5").Value = 0.23
$ws.Range("X5").Value = 0.115
$ws.Range("Y5").Value = 0.339
$ws.Range("AI6").Value = 0.219
$ws.Range("AU6").Value = 0.198
$ws.Range("BA6").Value = 1.628
$ws.Range("BG6").Value = 0.536
$ws.Range("BM6").Value = 0.636
$ws.Range("BP6").Value = 0.543
$ws.Range("BQ6").Value = 0.541
$ws.Range("E6").Value = 0.433
$ws.Range("N6").Value = 0.529
$ws.Range("W6").Value = 0.225
$ws.Range("AI7").Value = 0.231
$ws.Range("AU7").Value = 0.25
$ws.Range("BA7").Value = 1.477
$ws.Range("BG7").Value = 0.458
$ws.Range("BM7").Value = 0.604
$ws.Range("BP7").Value = 0.492
$ws.Range("BQ7").Value = 0.489
$ws.Range("E7").Value = 0.471
$ws.Range("N7").Value = 0.649
$ws.Range("W7").Value = 0.228
$ws.Range("AI8").Value = 0.22
$ws.Range("AJ8").Value = 0.093
$ws.Range("AK8").Value = 0.305
$ws.Range("AU8").Value = 0.244
$ws.Range("AV8").Value = 0.075
$ws.Range("AW8").Value = 0.275
$ws.Range("BA8").Value = 1.751
$ws.Range("BB8").Value = 0.134
$ws.Range("BC8").Value = 0.365
$ws.Range("BG8").Value = 0.577
$ws.Range("BH8").Value = 0.1
$ws.Range("BI8").Value = 0.317
$ws.Range("BM8").Value = 0.716
$ws.Range("BN8").Value = 0.063
$ws.Range("BO8").Value = 0.251
$ws.Range("BP8").Value = 0.584
$ws.Range("BQ8").Value = 0.588
$ws.Range("E8").Value = 0.525
$ws.Range("F8").Value = 0.129
$ws.Range("G8").Value = 0.359
$ws.Range("N8").Value = 0.761
$ws.Range("O8").Value = 0.062
$ws.Range("P8").Value = 0.248
$ws.Range("W8").Value = 0.222
$ws.Range("X8").Value = 0.107
$ws.Range("Y8").Value = 0.327
$ws.Range("AI9").Value = 0.122
$ws.Range("AJ9").Value = 0.107
$ws.Range("AK9").Value = 0.328
$ws.Range("BA9").Value = 1.694
$ws.Range("BB9").Value = 0.245
$ws.Range("BC9").Value = 0.495
$ws.Range("BG9").Value = 0.612
$ws.Range("BH9").Value = 0.237
$ws.Range("BI9").Value = 0.487
$ws.Range("BM9").Value = 0.653
$ws.Range("BN9").Value = 0.227
$ws.Range("BO9").Value = 0.476
$ws.Range("BP9").Value = 0.5649999999999999
$ws.Range("BQ9").Value = 0.554
$ws.Range("E9").Value = 0.449
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("N9").Value = 0.633
$ws.Range("O9").Value = 0.232
$ws.Range("P9").Value = 0.482
$ws.Range("W9").Value = 0.122
$ws.Range("X9").Value = 0.107
$ws.Range("Y9").Value = 0.328
$ws.Range("AI10").Value = 0.245
$ws.Range("AJ10").Value = 0.185
$ws.Range("AK10").Value = 0.43
$ws.Range("AU10").Value = 0.224
$ws.Range("AV10").Value = 0.174
$ws.Range("AW10").Value = 0.417
$ws.Range("BA10").Value = 2.021
$ws.Range("BB10").Value = 0.25
$ws.Range("BC10").Value = 0.5
$ws.Range("BG10").Value = 0.653
$ws.Range("BH10").Value = 0.227
$ws.Range("BI10").Value = 0.476
$ws.Range("BM10").Value = 0.878
$ws.Range("BN10").Value = 0.107
$ws.Range("BO10").Value = 0.328
$ws.Range("BP10").Value = 0.674
$ws.Range("BQ10").Value = 0.6889999999999999
$ws.Range("E10").Value = 0.571
$ws.Range("F10").Value = 0.245
$ws.Range("G10").Value = 0.495
$ws.Range("N10").Value = 0.837
$ws.Range("O10").Value = 0.137
$ws.Range("P10").Value = 0.37
$ws.Range("W10").Value = 0.265
$ws.Range("X10").Value = 0.195
$ws.Range("Y10").Value = 0.441
$ws.Range("AI11").Value = 0.245
$ws.Range("AJ11").Value = 0.185
$ws.Range("AK11").Value = 0.43
$ws.Range("AU11").Value = 0.347
$ws.Range("AV11").Value = 0.227
$ws.Range("AW11").Value = 0.476
$ws.Range("BA11").Value = 2.021
$ws.Range("BB11").Value = 0.25
$ws.Range("BC11").Value = 0.5
$ws.Range("BG11").Value = 0.653
$ws.Range("BH11").Value = 0.227
$ws.Range("BI11").Value = 0.476
$ws.Range("BM11").Value = 0.878
$ws.Range("BN11").Value = 0.107
$ws.Range("BO11").Value = 0.328
$ws.Range("BP11").Value = 0.674
$ws.Range("BQ11").Value = 0.6889999999999999
$ws.Range("E11").Value = 0.592
$ws.Range("F11").Value = 0.242
$ws.Range("G11").Value = 0.491
$ws.Range("N11").Value = 0.878
$ws.Range("O11").Value = 0.107
$ws.Range("P11").Value = 0.328
$ws.Range("W11").Value = 0.265
$ws.Range("X11").Value = 0.195
$ws.Range("Y11").Value = 0.441
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.647
$ws.Range("AV12").Value = 1.758
$ws.Range("AW12").Value = 1.326
$ws.Range("BA12").Value = 3.619
$ws.Range("BB12").Value = 0.332
$ws.Range("BC12").Value = 0.576
$ws.Range("BG12").Value = 1.062
$ws.Range("BH12").Value = 0.059
$ws.Range("BI12").Value = 0.242
$ws.Range("BM12").Value = 1.349
$ws.Range("BN12").Value = 0.413
$ws.Range("BO12").Value = 0.643
$ws.Range("BP12").Value = 1.206
$ws.Range("BQ12").Value = 1.262
$ws.Range("E12").Value = 1.414
$ws.Range("F12").Value = 0.656
$ws.Range("G12").Value = 0.8100000000000001
$ws.Range("N12").Value = 1.667
$ws.Range("O12").Value = 1.644
$ws.Range("P12").Value = 1.282
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.498
$ws.Range("AV13").Value = 1.32
$ws.Range("AW13").Value = 1.149
$ws.Range("BA13").Value = 2.523
$ws.Range("BB13").Value = 0.319
$ws.Range("BC13").Value = 0.5649999999999999
$ws.Range("BG13").Value = 0.627
$ws.Range("BH13").Value = 0.08699999999999999
$ws.Range("BI13").Value = 0.294
$ws.Range("BM13").Value = 0.999
$ws.Range("BN13").Value = 0.37
$ws.Range("BO13").Value = 0.609
$ws.Range("BP13").Value = 0.841
$ws.Range("BQ13").Value = 0.789
$ws.Range("E13").Value = 1.75
$ws.Range("F13").Value = 0.915
$ws.Range("G13").Value = 0.957
$ws.Range("N13").Value = 2.345
$ws.Range("O13").Value = 1.188
$ws.Range("P13").Value = 1.09
$ws.Range("W13").Value = 1.11
$ws.Range("X13").Value = 0.175
$ws.Range("Y13").Value = 0.419
